$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) and "全部类型" (index 4) both received updated
# "想去人数" (F column) counts, and row 22 also got an updated
# "最低票价" (G column) value.

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 441
$ws1.Range("F4").Value = 153
$ws1.Range("F5").Value = 1832
$ws1.Range("F6").Value = 1437
$ws1.Range("F8").Value = 1710
$ws1.Range("F10").Value = 138
$ws1.Range("F11").Value = 641
$ws1.Range("F12").Value = 26
$ws1.Range("F16").Value = 70
$ws1.Range("F17").Value = 133
$ws1.Range("F20").Value = 61
$ws1.Range("F21").Value = 98
$ws1.Range("F22").Value = 4414
$ws1.Range("G22").Value = 65
$ws1.Range("F23").Value = 33
$ws1.Range("F24").Value = 800
$ws1.Range("F25").Value = 94
$ws1.Range("F26").Value = 2122
$ws1.Range("F27").Value = 65
$ws1.Range("F28").Value = 1998

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 441
$ws4.Range("F4").Value = 153
$ws4.Range("F5").Value = 1832
$ws4.Range("F6").Value = 1437
$ws4.Range("F8").Value = 1710
$ws4.Range("F10").Value = 138
$ws4.Range("F11").Value = 641
$ws4.Range("F12").Value = 26
$ws4.Range("F16").Value = 70
$ws4.Range("F17").Value = 133
$ws4.Range("F20").Value = 61
$ws4.Range("F21").Value = 98
$ws4.Range("F22").Value = 4414
$ws4.Range("G22").Value = 65
$ws4.Range("F24").Value = 33
$ws4.Range("F26").Value = 800
$ws4.Range("F27").Value = 94
$ws4.Range("F28").Value = 2122
$ws4.Range("F29").Value = 65
$ws4.Range("F30").Value = 1998
